# Updated cryptos list on Wed May 29 22:58:11 UTC 2024 with GitHub Actions
# Refreshes price/volume columns for each coin row and swaps the
# Kaspa/Mantle rows (rows 39-40) to reflect the new ranking order.
# NumberFormat is forced to text ("@") before writing any price value that
# Excel would otherwise auto-parse as a number, then the cell style is
# reset back to Normal so formatting/appearance is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.673.44'
$ws.Range('E2').Value = '  -1.30%  '
$ws.Range('D3').Value = '3.779.36'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.21'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '3.776.44'
$ws.Range('E7').Value = '  -2.20%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('E12').Value = '  -1.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000278'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.49'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('D15').Value = '4.412.50'
$ws.Range('E15').Value = '  -2.13%  '
$ws.Range('D16').Value = '3.779.17'
$ws.Range('E16').Value = '  -2.21%  '
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').Value = '67.645.28'
$ws.Range('E18').Value = '  -1.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.83%  '
$ws.Range('E20').Value = '  +0.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.51'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '468.64'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.58%  '
$ws.Range('E23').Value = '  -2.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000148'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.82'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('E26').Value = '  -1.56%  '
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.29'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.43%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('E30').Value = '  -1.69%  '
$ws.Range('D31').Value = '3.931.16'
$ws.Range('E31').Value = '  -2.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.61'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '30.53'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.16%  '
$ws.Range('E34').Value = '  -3.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.17%  '
$ws.Range('D36').Value = '3.741.80'
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.76'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.27%  '
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.01'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.60%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.138'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.80'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('E43').Value = '  -1.08%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.94'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '45.86'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '396.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000270'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '140.86'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '39.22'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.20%  '
